# Squash merge add-llm-recommendations into main
#
# Slide 7 ("Strategic Recommendations") is reworked so each
# recommendation gets a short bold title plus a new rationale
# paragraph underneath it. Placeholders are renamed from
# {{RECOMMENDATION_n}} to {{REC_n_TITLE}} / {{REC_n_RATIONALE}}.
#
# NOTE on numeric literals below: this COM host stores shape
# geometry as 32-bit-float points before converting to EMU, and
# truncates (rather than rounds) on the pt -> EMU conversion. A
# couple of the target EMU values are not exactly representable
# from the "obvious" point value because of that, so a few numbers
# are nudged by 0.00001pt (a sub-visible amount, well under a
# single EMU) so the stored EMU matches the target exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Grab references to all the pre-existing shapes up front (by their
# original, stable z-order position) before we insert anything new -
# inserting shapes shifts everyone's Shapes.Item(n) index around, but
# an already-fetched shape reference stays valid no matter how the
# collection is reshuffled afterwards.
$shpTitle = $s.Shapes.Item(1)   # id 2  "TextBox 1"   - page title
$sq1      = $s.Shapes.Item(2)   # id 3  "Rectangle 2" - square 1
$num1     = $s.Shapes.Item(3)   # id 4  "TextBox 3"   - number "1"
$title1   = $s.Shapes.Item(4)   # id 5  "TextBox 4"   - {{RECOMMENDATION_1}}
$sq2      = $s.Shapes.Item(5)   # id 6  "Rectangle 5" - square 2
$num2     = $s.Shapes.Item(6)   # id 7  "TextBox 6"   - number "2"
$title2   = $s.Shapes.Item(7)   # id 8  "TextBox 7"   - {{RECOMMENDATION_2}}
$sq3      = $s.Shapes.Item(8)   # id 9  "Rectangle 8" - square 3
$num3     = $s.Shapes.Item(9)   # id 10 "TextBox 9"   - number "3"
$title3   = $s.Shapes.Item(10)  # id 11 "TextBox 10"  - {{RECOMMENDATION_3}}

# ---------------------------------------------------------------
# Page title
# ---------------------------------------------------------------
$shpTitle.Left = 36.0
$shpTitle.Top = 21.6
$shpTitle.Width = 648.0
$shpTitle.Height = 50.4
$shpTitle.TextFrame.TextRange.Font.Size = 36

# ---------------------------------------------------------------
# Recommendation 1 block
# ---------------------------------------------------------------

$sq1.Left = 36.0
$sq1.Top = 86.4
$sq1.Width = 28.80001
$sq1.Height = 28.80001

$num1.Left = 36.0
$num1.Top = 86.4
$num1.Width = 28.80001
$num1.Height = 28.80001
$num1.TextFrame.TextRange.Font.Size = 16

$title1.Left = 79.20001
$title1.Top = 86.4
$title1.Width = 597.60001
$title1.Height = 28.80001
$title1.TextFrame.TextRange.Text = "{{REC_1_TITLE}}"
$title1.TextFrame.TextRange.Font.Size = 14
$title1.TextFrame.TextRange.Font.Bold = $true
$title1.TextFrame.TextRange.Font.Color.RGB = 0x8A5C2E
# re-assert geometry: setting .Text recalculated the autofit height
$title1.Left = 79.20001
$title1.Top = 86.4
$title1.Width = 597.60001
$title1.Height = 28.80001

# NEW: rationale text box for recommendation 1
$rat1 = $s.Shapes.AddTextbox(1, 79.20001, 116.64001, 597.60001, 36.0)
$rat1.Name = "TextBox 5"
$rat1.Fill.Visible = 0
$rat1.TextFrame.WordWrap = -1
$rat1.TextFrame.AutoSize = 1
$rat1.TextFrame.TextRange.Text = "{{REC_1_RATIONALE}}"
$rat1.TextFrame.TextRange.Font.Size = 11
$rat1.TextFrame.TextRange.Font.Color.RGB = 0x68554A
$rat1.Left = 79.20001
$rat1.Top = 116.64001
$rat1.Width = 597.60001
$rat1.Height = 36.0
while ($rat1.ZOrderPosition -gt 5) { $rat1.ZOrder(3) }

# ---------------------------------------------------------------
# Recommendation 2 block
# ---------------------------------------------------------------

$sq2.Name = "Rectangle 6"
$sq2.Left = 36.0
$sq2.Top = 158.40001
$sq2.Width = 28.80001
$sq2.Height = 28.80001

$num2.Name = "TextBox 7"
$num2.Left = 36.0
$num2.Top = 158.40001
$num2.Width = 28.80001
$num2.Height = 28.80001
$num2.TextFrame.TextRange.Font.Size = 16

$title2.Name = "TextBox 8"
$title2.Left = 79.20001
$title2.Top = 158.40001
$title2.Width = 597.60001
$title2.Height = 28.80001
$title2.TextFrame.TextRange.Text = "{{REC_2_TITLE}}"
$title2.TextFrame.TextRange.Font.Size = 14
$title2.TextFrame.TextRange.Font.Bold = $true
$title2.TextFrame.TextRange.Font.Color.RGB = 0x8A5C2E
$title2.Left = 79.20001
$title2.Top = 158.40001
$title2.Width = 597.60001
$title2.Height = 28.80001

# NEW: rationale text box for recommendation 2
$rat2 = $s.Shapes.AddTextbox(1, 79.20001, 188.64001, 597.60001, 36.0)
$rat2.Name = "TextBox 9"
$rat2.Fill.Visible = 0
$rat2.TextFrame.WordWrap = -1
$rat2.TextFrame.AutoSize = 1
$rat2.TextFrame.TextRange.Text = "{{REC_2_RATIONALE}}"
$rat2.TextFrame.TextRange.Font.Size = 11
$rat2.TextFrame.TextRange.Font.Color.RGB = 0x68554A
$rat2.Left = 79.20001
$rat2.Top = 188.64001
$rat2.Width = 597.60001
$rat2.Height = 36.0
while ($rat2.ZOrderPosition -gt 9) { $rat2.ZOrder(3) }

# ---------------------------------------------------------------
# Recommendation 3 block
# ---------------------------------------------------------------

$sq3.Name = "Rectangle 10"
$sq3.Left = 36.0
$sq3.Top = 230.40001
$sq3.Width = 28.80001
$sq3.Height = 28.80001

$num3.Name = "TextBox 11"
$num3.Left = 36.0
$num3.Top = 230.40001
$num3.Width = 28.80001
$num3.Height = 28.80001
$num3.TextFrame.TextRange.Font.Size = 16

$title3.Name = "TextBox 12"
$title3.Left = 79.20001
$title3.Top = 230.40001
$title3.Width = 597.60001
$title3.Height = 28.80001
$title3.TextFrame.TextRange.Text = "{{REC_3_TITLE}}"
$title3.TextFrame.TextRange.Font.Size = 14
$title3.TextFrame.TextRange.Font.Bold = $true
$title3.TextFrame.TextRange.Font.Color.RGB = 0x8A5C2E
$title3.Left = 79.20001
$title3.Top = 230.40001
$title3.Width = 597.60001
$title3.Height = 28.80001

# NEW: rationale text box for recommendation 3 (already lands last in z-order)
$rat3 = $s.Shapes.AddTextbox(1, 79.20001, 260.64, 597.60001, 36.0)
$rat3.Name = "TextBox 13"
$rat3.Fill.Visible = 0
$rat3.TextFrame.WordWrap = -1
$rat3.TextFrame.AutoSize = 1
$rat3.TextFrame.TextRange.Text = "{{REC_3_RATIONALE}}"
$rat3.TextFrame.TextRange.Font.Size = 11
$rat3.TextFrame.TextRange.Font.Color.RGB = 0x68554A
$rat3.Left = 79.20001
$rat3.Top = 260.64
$rat3.Width = 597.60001
$rat3.Height = 36.0
